{"js": "// Fix the description of generating the HTML report: the experiment report\n// section should describe fetching data about a run (not \"the experiment\"),\n// the application \"fetches data for all runs\" (instead of just \"creates\"),\n// and the html file is created \"to print the data to,\" (instead of\n// \"and prints the data\").\nconst body = context.document.body;\nconst target =\n  \"fetching data about the experiment except the application creates an html file and prints the data along with\";\nconst replacement =\n  \"fetching data about a run except the application fetches data for all runs and creates an html file to print the data to, along with\";\n\nconst results = body.search(target, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Fix the description of generating the HTML report: the experiment report\n# section should describe fetching data about a run (not \"the experiment\"),\n# the application \"fetches data for all runs\" (instead of just \"creates\"),\n# and the html file is created \"to print the data to,\" (instead of\n# \"and prints the data\").\n$d = $word.ActiveDocument\n\n$searchText = \"fetching data about the experiment except the application creates an html file and prints the data along with\"\n$replaceText = \"fetching data about a run except the application fetches data for all runs and creates an html file to print the data to, along with\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $searchText\n$find.Replacement.Text = $replaceText\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
